$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update column widths (B=47, D=28, H=13 in raw XML terms)
$ws.Columns.Item(2).ColumnWidth = 46.17
$ws.Columns.Item(4).ColumnWidth = 27.17
$ws.Columns.Item(8).ColumnWidth = 12.17

# Remove all existing hyperlinks before rewriting URLs so relationship targets refresh
$ws.Hyperlinks.Delete()

# Row 2
$ws.Cells.Item(2, 1).Value = '2025-09-22 06:27:57'
$ws.Cells.Item(2, 2).Value = '【せどり×ツール製作】APIを使用したせどりツールを製作できるエンジニアさんを募集します♪'
$ws.Cells.Item(2, 3).Value = 'システム開発'
$ws.Cells.Item(2, 4).Value = '20,000 円 ~ 50,000 円 / 固定'
$ws.Cells.Item(2, 5).Value = '期限情報なし'
$ws.Cells.Item(2, 6).Value = 'https://www.lancers.jp/work/detail/5217096'
$ws.Cells.Item(2, 7).Value = 243
$ws.Cells.Item(2, 8).Value = '🔥API ◆ツール'

# Row 3
$ws.Cells.Item(3, 1).Value = '2025-09-22 06:27:57'
$ws.Cells.Item(3, 2).Value = '【相談希望】在庫管理・出品補助ツールの開発に関するZoom面談依頼'
$ws.Cells.Item(3, 3).Value = 'システム開発'
$ws.Cells.Item(3, 4).Value = '20,000 円 ~ 50,000 円 / 固定'
$ws.Cells.Item(3, 5).Value = '期限情報なし'
$ws.Cells.Item(3, 6).Value = 'https://www.lancers.jp/work/detail/5398112'
$ws.Cells.Item(3, 7).Value = 158
$ws.Cells.Item(3, 8).Value = '◆ツール,開発 ◇管理'

# Row 4
$ws.Cells.Item(4, 1).Value = '2025-09-22 06:27:57'
$ws.Cells.Item(4, 2).Value = '【急募】Slack自動リアクションツール開発依頼'
$ws.Cells.Item(4, 3).Value = 'システム開発'
$ws.Cells.Item(4, 4).Value = '5,000 円 ~ 10,000 円 / 固定'
$ws.Cells.Item(4, 5).Value = '期限情報なし'
$ws.Cells.Item(4, 6).Value = 'https://www.lancers.jp/work/detail/5398193'
$ws.Cells.Item(4, 7).Value = 120
$ws.Cells.Item(4, 8).Value = '◆ツール,開発'

# Row 5
$ws.Cells.Item(5, 1).Value = '2025-09-22 06:27:57'
$ws.Cells.Item(5, 2).Value = '【急募】MT4特定口座の取引を子口座に反映するシステム開発'
$ws.Cells.Item(5, 3).Value = 'システム開発'
$ws.Cells.Item(5, 4).Value = '100,000 円 ~ 200,000 円 / 固定'
$ws.Cells.Item(5, 5).Value = '期限情報なし'
$ws.Cells.Item(5, 6).Value = 'https://www.lancers.jp/work/detail/5398203'
$ws.Cells.Item(5, 7).Value = 118
$ws.Cells.Item(5, 8).Value = '◆開発,システム開発'

# Row 6
$ws.Cells.Item(6, 1).Value = '2025-09-22 06:27:57'
$ws.Cells.Item(6, 2).Value = '【急募】自己分析アプリのバックエンド開発アドバイザリー募集'
$ws.Cells.Item(6, 3).Value = 'システム開発'
$ws.Cells.Item(6, 4).Value = '50,000 円 ~ 100,000 円 / 固定'
$ws.Cells.Item(6, 5).Value = '期限情報なし'
$ws.Cells.Item(6, 6).Value = 'https://www.lancers.jp/work/detail/5397930'
$ws.Cells.Item(6, 7).Value = 93
$ws.Cells.Item(6, 8).Value = '◆開発 ◇アプリ'

# Row 7
$ws.Cells.Item(7, 1).Value = '2025-09-22 06:27:57'
$ws.Cells.Item(7, 2).Value = '【GAS開発者募集】Amazon広告管理SaaSのMVP開発'
$ws.Cells.Item(7, 3).Value = 'システム開発'
$ws.Cells.Item(7, 4).Value = '20,000 円 ~ 50,000 円 / 固定'
$ws.Cells.Item(7, 5).Value = '期限情報なし'
$ws.Cells.Item(7, 6).Value = 'https://www.lancers.jp/work/detail/5397812'
$ws.Cells.Item(7, 7).Value = 88
$ws.Cells.Item(7, 8).Value = '◆開発 ◇管理'

# Row 8
$ws.Cells.Item(8, 1).Value = '2025-09-22 06:27:57'
$ws.Cells.Item(8, 2).Value = '【急募】iOSアプリのAdMobメディエーション入札接続とeCPM改善'
$ws.Cells.Item(8, 3).Value = 'システム開発'
$ws.Cells.Item(8, 4).Value = '50,000 円 ~ 100,000 円 / 固定'
$ws.Cells.Item(8, 5).Value = '期限情報なし'
$ws.Cells.Item(8, 6).Value = 'https://www.lancers.jp/work/detail/5398081'
$ws.Cells.Item(8, 7).Value = 38
$ws.Cells.Item(8, 8).Value = '◇アプリ'

# Row 9
$ws.Cells.Item(9, 1).Value = '2025-09-22 06:27:57'
$ws.Cells.Item(9, 2).Value = '【Braze経験者募集】CRM/マーケティングオートメーション支援(中級者以上)'
$ws.Cells.Item(9, 3).Value = 'システム開発'
$ws.Cells.Item(9, 4).Value = '300,000 円 ~ 500,000 円 / 固定'
$ws.Cells.Item(9, 5).Value = '期限情報なし'
$ws.Cells.Item(9, 6).Value = 'https://www.lancers.jp/work/detail/5398071'
$ws.Cells.Item(9, 7).Value = 25
$ws.Cells.Item(9, 8).ClearContents()

# Row 10
$ws.Cells.Item(10, 1).Value = '2025-09-22 06:27:57'
$ws.Cells.Item(10, 2).Value = '【Braze経験者募集】CRM/マーケティングオートメーション支援(中級者以上)'
$ws.Cells.Item(10, 3).Value = 'システム開発'
$ws.Cells.Item(10, 4).Value = '300,000 円 ~ 500,000 円 / 固定'
$ws.Cells.Item(10, 5).Value = '期限情報なし'
$ws.Cells.Item(10, 6).Value = 'https://www.lancers.jp/work/detail/5398062'
$ws.Cells.Item(10, 7).Value = 25
$ws.Cells.Item(10, 8).ClearContents()

# Row 11
$ws.Cells.Item(11, 1).Value = '2025-09-22 06:27:57'
$ws.Cells.Item(11, 2).Value = 'データセンター向けサーバー・ルーター設置作業'
$ws.Cells.Item(11, 3).Value = 'システム開発'
$ws.Cells.Item(11, 4).Value = '100,000 円 ~ 200,000 円 / 固定'
$ws.Cells.Item(11, 5).Value = '期限情報なし'
$ws.Cells.Item(11, 6).Value = 'https://www.lancers.jp/work/detail/5397887'
$ws.Cells.Item(11, 7).Value = 18
$ws.Cells.Item(11, 8).ClearContents()

# Row 12
$ws.Cells.Item(12, 1).Value = '2025-09-22 06:27:57'
$ws.Cells.Item(12, 2).Value = 'Excelやスプレッドシートでのデータシュミレーション クエリ(query)や関数利用'
$ws.Cells.Item(12, 3).Value = 'システム開発'
$ws.Cells.Item(12, 4).Value = '5,000 円 ~ 10,000 円 / 固定'
$ws.Cells.Item(12, 5).Value = '期限情報なし'
$ws.Cells.Item(12, 6).Value = 'https://www.lancers.jp/work/detail/5397980'
$ws.Cells.Item(12, 7).Value = 10
$ws.Cells.Item(12, 8).ClearContents()

# Row 13
$ws.Cells.Item(13, 1).Value = '2025-09-22 06:27:57'
$ws.Cells.Item(13, 2).Value = '【中小企業支援】債務超過・赤字経営解消の診断依頼'
$ws.Cells.Item(13, 3).Value = 'システム開発'
$ws.Cells.Item(13, 4).Value = '5,000 円 ~ 10,000 円 / 固定'
$ws.Cells.Item(13, 5).Value = '期限情報なし'
$ws.Cells.Item(13, 6).Value = 'https://www.lancers.jp/work/detail/5397962'
$ws.Cells.Item(13, 7).Value = 10
$ws.Cells.Item(13, 8).ClearContents()

# Row 14
$ws.Cells.Item(14, 1).Value = '2025-09-22 06:27:57'
$ws.Cells.Item(14, 2).Value = 'Geminiで旅行のしおりのHTMLを生成するプロンプトの作成'
$ws.Cells.Item(14, 3).Value = 'システム開発'
$ws.Cells.Item(14, 4).Value = '10,000 円 ~ 20,000 円 / 固定'
$ws.Cells.Item(14, 5).Value = '期限情報なし'
$ws.Cells.Item(14, 6).Value = 'https://www.lancers.jp/work/detail/5397817'
$ws.Cells.Item(14, 7).Value = 10
$ws.Cells.Item(14, 8).ClearContents()

# Re-add hyperlinks for F2:F14 in order, then restyle as Hyperlink
$ws.Hyperlinks.Add($ws.Cells.Item(2, 6), 'https://www.lancers.jp/work/detail/5217096')
$ws.Hyperlinks.Add($ws.Cells.Item(3, 6), 'https://www.lancers.jp/work/detail/5398112')
$ws.Hyperlinks.Add($ws.Cells.Item(4, 6), 'https://www.lancers.jp/work/detail/5398193')
$ws.Hyperlinks.Add($ws.Cells.Item(5, 6), 'https://www.lancers.jp/work/detail/5398203')
$ws.Hyperlinks.Add($ws.Cells.Item(6, 6), 'https://www.lancers.jp/work/detail/5397930')
$ws.Hyperlinks.Add($ws.Cells.Item(7, 6), 'https://www.lancers.jp/work/detail/5397812')
$ws.Hyperlinks.Add($ws.Cells.Item(8, 6), 'https://www.lancers.jp/work/detail/5398081')
$ws.Hyperlinks.Add($ws.Cells.Item(9, 6), 'https://www.lancers.jp/work/detail/5398071')
$ws.Hyperlinks.Add($ws.Cells.Item(10, 6), 'https://www.lancers.jp/work/detail/5398062')
$ws.Hyperlinks.Add($ws.Cells.Item(11, 6), 'https://www.lancers.jp/work/detail/5397887')
$ws.Hyperlinks.Add($ws.Cells.Item(12, 6), 'https://www.lancers.jp/work/detail/5397980')
$ws.Hyperlinks.Add($ws.Cells.Item(13, 6), 'https://www.lancers.jp/work/detail/5397962')
$ws.Hyperlinks.Add($ws.Cells.Item(14, 6), 'https://www.lancers.jp/work/detail/5397817')
$ws.Range("F2:F14").Style = "Hyperlink"

Write-Host "done"